$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws1.Range("H70").Value = 11919.667
$ws1.Range("I70").Value = 3519.2
$ws1.Range("J70").Value = 16119.9
$ws1.Range("K70").Value = 10557.6
$ws1.Range("L70").Value = 48359.7
$ws1.Range("M70").Value = -10287.6
$ws1.Range("N70").Value = -48899.7

$ws1.Range("H73").Value = 11919.667
$ws1.Range("I73").Value = 3519.2
$ws1.Range("J73").Value = 16119.9
$ws1.Range("K73").Value = 10557.6
$ws1.Range("L73").Value = 48359.7
$ws1.Range("M73").Value = -9621.599999999999
$ws1.Range("N73").Value = -50231.7

$ws1.Range("H98").Value = 1836275.2
$ws1.Range("I98").Value = 1977752
$ws1.Range("K98").Value = 1977752
$ws1.Range("M98").Value = -1976254

$ws1.Range("H107").Value = 9805772
$ws1.Range("I107").Value = 6945794.5
$ws1.Range("K107").Value = 6945794.5
$ws1.Range("M107").Value = -6943874.5

$ws1.Range("H122").Value = 1836275.2
$ws1.Range("I122").Value = 1977752
$ws1.Range("K122").Value = 5933256
$ws1.Range("M122").Value = -5930806

$ws1.Range("H131").Value = 37383.8
$ws1.Range("I131").Value = 2365.0557
$ws1.Range("K131").Value = 7095.1671
$ws1.Range("M131").Value = -2055.1671

$ws1.Range("H132").Value = 1858.0238
$ws1.Range("I132").Value = 1835.4242
$ws1.Range("J132").Value = 1940.8889
$ws1.Range("K132").Value = 5506.2726
$ws1.Range("L132").Value = 5822.6667
$ws1.Range("M132").Value = -2976.2726
$ws1.Range("N132").Value = -10882.6667

$ws1.Range("H138").Value = 2328.57
$ws1.Range("I138").Value = 1792.9143
$ws1.Range("J138").Value = 2617
$ws1.Range("K138").Value = 5378.742899999999
$ws1.Range("L138").Value = 7851
$ws1.Range("M138").Value = -238.7428999999993
$ws1.Range("N138").Value = -18131

$ws2 = $wb.Worksheets.Item("ARM")
$ws2.Range("H32").Value = 13955.181
$ws2.Range("I32").Value = 12646.597
$ws2.Range("J32").Value = 32602.5
$ws2.Range("K32").Value = 12646.597
$ws2.Range("L32").Value = 32602.5
$ws2.Range("M32").Value = -12359.597
$ws2.Range("N32").Value = -33176.5

$ws2.Range("H61").Value = 1640.3392
$ws2.Range("I61").Value = 1268.0588
$ws2.Range("J61").Value = 5437.6
$ws2.Range("K61").Value = 1268.0588
$ws2.Range("L61").Value = 5437.6
$ws2.Range("M61").Value = -1056.0588
$ws2.Range("N61").Value = -5861.6

$ws2.Range("H132").Value = 3360.8708
$ws2.Range("I132").Value = 2511.4333
$ws2.Range("J132").Value = 28844
$ws2.Range("K132").Value = 7534.2999
$ws2.Range("L132").Value = 86532
$ws2.Range("M132").Value = -5004.2999
$ws2.Range("N132").Value = -91592

$ws2.Range("H136").Value = 1640.3392
$ws2.Range("I136").Value = 1268.0588
$ws2.Range("J136").Value = 5437.6
$ws2.Range("K136").Value = 3804.1764
$ws2.Range("L136").Value = 16312.8
$ws2.Range("M136").Value = -1254.1764
$ws2.Range("N136").Value = -21412.8

$ws3 = $wb.Worksheets.Item("BSM")
$ws3.Range("H86").Value = 1720.8667
$ws3.Range("I86").Value = 1630.5217
$ws3.Range("J86").Value = 2017.7142
$ws3.Range("K86").Value = 1630.5217
$ws3.Range("L86").Value = 2017.7142
$ws3.Range("M86").Value = -507.5217
$ws3.Range("N86").Value = -4263.7142

$ws3.Range("H89").Value = 1720.8667
$ws3.Range("I89").Value = 1630.5217
$ws3.Range("J89").Value = 2017.7142
$ws3.Range("K89").Value = 8152.6085
$ws3.Range("L89").Value = 10088.571
$ws3.Range("M89").Value = -2536.6085
$ws3.Range("N89").Value = -21320.571

$ws3.Range("H134").Value = 1967.0247
$ws3.Range("I134").Value = 1561.5
$ws3.Range("K134").Value = 4684.5
$ws3.Range("M134").Value = -2149.5

$ws4 = $wb.Worksheets.Item("CRP")
$ws4.Range("H31").Value = 30388.896
$ws4.Range("I31").Value = 2546.5293
$ws4.Range("K31").Value = 2546.5293
$ws4.Range("M31").Value = -2251.5293

$ws4.Range("H34").Value = 30388.896
$ws4.Range("I34").Value = 2546.5293
$ws4.Range("K34").Value = 2546.5293
$ws4.Range("M34").Value = -2344.5293

$ws4.Range("H105").Value = 4750.0527
$ws4.Range("I105").Value = 1906.5294
$ws4.Range("K105").Value = 1906.5294
$ws4.Range("M105").Value = -159.5293999999999

$ws4.Range("H107").Value = 4368.5527
$ws4.Range("I107").Value = 489.77777
$ws4.Range("J107").Value = 7859.45
$ws4.Range("K107").Value = 489.77777
$ws4.Range("L107").Value = 7859.45
$ws4.Range("M107").Value = 1430.22223
$ws4.Range("N107").Value = -11699.45

$ws4.Range("H132").Value = 3391.2593
$ws4.Range("I132").Value = 2960.889
$ws4.Range("K132").Value = 8882.667000000001
$ws4.Range("M132").Value = -6352.667000000001

$ws5 = $wb.Worksheets.Item("CUL")
$ws5.Range("H113").Value = 2166.9583
$ws5.Range("J113").Value = 2179.8635
$ws5.Range("L113").Value = 6539.5905
$ws5.Range("N113").Value = -10879.5905

$ws5.Range("H131").Value = 8334975.5
$ws5.Range("I131").Value = 50000576
$ws5.Range("J131").Value = 1855.62
$ws5.Range("K131").Value = 150001728
$ws5.Range("L131").Value = 5566.86
$ws5.Range("M131").Value = -149996688
$ws5.Range("N131").Value = -15646.86

$ws6 = $wb.Worksheets.Item("GSM")
$ws6.Range("H102").Value = 1879.711
$ws6.Range("I102").Value = 1876.9773
$ws6.Range("K102").Value = 1876.9773
$ws6.Range("M102").Value = -254.9773

$ws6.Range("H132").Value = 28376.281
$ws6.Range("I132").Value = 29717.8
$ws6.Range("J132").Value = 16638
$ws6.Range("K132").Value = 89153.39999999999
$ws6.Range("L132").Value = 49914
$ws6.Range("M132").Value = -86623.39999999999
$ws6.Range("N132").Value = -54974

$ws7 = $wb.Worksheets.Item("LTW")
$ws7.Range("H16").Value = 2288.2104
$ws7.Range("I16").Value = 1424.4688
$ws7.Range("J16").Value = 6894.8335
$ws7.Range("K16").Value = 1424.4688
$ws7.Range("L16").Value = 6894.8335
$ws7.Range("M16").Value = -1254.4688
$ws7.Range("N16").Value = -7234.8335

$ws7.Range("H22").Value = 2642.875
$ws7.Range("I22").Value = 1866
$ws7.Range("J22").Value = 3109
$ws7.Range("K22").Value = 1866
$ws7.Range("L22").Value = 3109
$ws7.Range("M22").Value = -1571
$ws7.Range("N22").Value = -3699

$ws7.Range("H27").Value = 2642.875
$ws7.Range("I27").Value = 1866
$ws7.Range("J27").Value = 3109
$ws7.Range("K27").Value = 1866
$ws7.Range("L27").Value = 3109
$ws7.Range("M27").Value = -1759
$ws7.Range("N27").Value = -3323

$ws7.Range("H132").Value = 6055.5776
$ws7.Range("I132").Value = 1557.9524
$ws7.Range("J132").Value = 9991
$ws7.Range("K132").Value = 4673.857199999999
$ws7.Range("L132").Value = 29973
$ws7.Range("M132").Value = -2143.857199999999
$ws7.Range("N132").Value = -35033

$ws8 = $wb.Worksheets.Item("WVR")
$ws8.Range("H62").Value = 12995
$ws8.Range("I62").Value = 0
$ws8.Range("J62").Value = 12995
$ws8.Range("K62").Value = 0
$ws8.Range("L62").Value = 12995
$ws8.Range("M62").ClearContents()
$ws8.Range("N62").Value = -14243

$ws8.Range("H65").Value = 12995
$ws8.Range("I65").Value = 0
$ws8.Range("J65").Value = 12995
$ws8.Range("K65").Value = 0
$ws8.Range("L65").Value = 64975
$ws8.Range("M65").ClearContents()
$ws8.Range("N65").Value = -71215

$ws8.Range("H107").Value = 456.3913
$ws8.Range("J107").Value = 741
$ws8.Range("L107").Value = 2223
$ws8.Range("N107").Value = -6063

$ws8.Range("H132").Value = 1648.7234
$ws8.Range("I132").Value = 907.88464
$ws8.Range("K132").Value = 2723.65392
$ws8.Range("M132").Value = -193.6539199999997
